$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "loading_percent" results for Case_4_31 (380 kV case).
# Each entry: @(rowNumber, @(@(colNumber, newValue), ...))
$data = @(
    @(2, @(@(2,14.86880019082857),@(3,10.00077242484913),@(4,5.973662031265871),@(5,11.49678251927395),@(7,26.56228231048298),@(8,13.56233498213508),@(9,20.58843904220758),@(12,9.726813759130513),@(14,16.81341884539283),@(15,20.38161832807902))),
    @(3, @(@(2,14.32207975470207),@(3,9.789685452692897),@(4,5.852700298432294),@(5,11.53498044091433),@(7,26.46717264402335),@(8,13.59807548138736),@(9,20.67889318057254),@(12,9.699109856417412),@(14,16.85506813843295),@(15,20.41267717606436))),
    @(4, @(@(2,13.97701550734586),@(3,9.656692340611436),@(4,5.778920700046474),@(5,11.56038330156764),@(7,26.41897708514221),@(8,13.62265925009512),@(9,20.73913748935086),@(12,9.683786179357698),@(14,16.88238532226372),@(15,20.43716364719482))),
    @(5, @(@(2,13.83425110134348),@(3,9.60169351449778),@(4,5.749025591324566),@(5,11.57122517646601),@(7,26.40191724767864),@(8,13.63333993923401),@(9,20.76486807327877),@(12,9.67797033434835),@(14,16.89395664976177),@(15,20.4485007522205))),
    @(6, @(@(2,13.81042200149882),@(3,9.592513953441427),@(4,5.744073226861144),@(5,11.57305505252047),@(7,26.39924066566707),@(8,13.63515344910354),@(9,20.76921184286478),@(12,9.677030634246822),@(14,16.8959046172347),@(15,20.45046522338398))),
    @(7, @(@(2,13.9750985406363),@(3,9.655953793354104),@(4,5.778516770433387),@(5,11.56052753505545),@(7,26.41873654646438),@(8,13.62280061210695),@(9,20.73947972449339),@(12,9.683706003407099),@(14,16.88253959734452),@(15,20.43731104692419))),
    @(8, @(@(2,14.68235608290466),@(3,9.928720229764041),@(4,5.931880501638131),@(5,11.50954853119672),@(7,26.52738121710543),@(8,13.57411001869289),@(9,20.6186495271386),@(12,9.716914246070504),@(14,16.82741791341908),@(15,20.39120147360258))),
    @(9, @(@(2,15.98675503912645),@(3,10.43477703056634),@(4,6.23449588386173),@(5,11.42505058466735),@(7,26.82052979456945),@(8,13.49960870327734),@(9,20.41915048715371),@(12,9.795201023790058),@(14,16.73313312234053),@(15,20.34387943478038))),
    @(10, @(@(2,16.88530637234834),@(3,10.78653768304301),@(4,6.455337781299679),@(5,11.37240960157923),@(7,27.08330968248906),@(8,13.45771699665983),@(9,20.29556599785491),@(12,9.860429767444128),@(14,16.67223689404519),@(15,20.3355200446265))),
    @(11, @(@(2,17.27947035319913),@(3,10.9417459148471),@(4,6.554958354641228),@(5,11.35051318889866),@(7,27.21276238921976),@(8,13.44145845188706),@(9,20.24436954623339),@(12,9.891703884235504),@(14,16.64634325482855),@(15,20.33746718416446))),
    @(12, @(@(2,17.42651860154309),@(3,10.99979237725337),@(4,6.592520383998105),@(5,11.34251656394632),@(7,27.26317032184013),@(8,13.4357048229868),@(9,20.22570789438939),@(12,9.90376997626988),@(14,16.63679740188503),@(15,20.33903147987271))),
    @(13, @(@(2,17.39494922762208),@(3,10.98732393303785),@(4,6.584438568592526),@(5,11.34422565293899),@(7,27.25225303144552),@(8,13.43692602967965),@(9,20.2296947069495),@(12,9.901161512481128),@(14,16.63884173881201),@(15,20.33865780693937))),
    @(14, @(@(2,17.29161305923935),@(3,10.94653620906688),@(4,6.55805208158486),@(5,11.34984938667872),@(7,27.21688191451524),@(8,13.44097701409884),@(9,20.24281968663351),@(12,9.892692135545769),@(14,16.64555271386404),@(15,20.33757930797521))),
    @(15, @(@(2,17.22802531087024),@(3,10.92145678522052),@(4,6.541867303835175),@(5,11.35333252025633),@(7,27.19539551518435),@(8,13.44351087902692),@(9,20.25095366623423),@(12,9.887533260086538),@(14,16.64969715969677),@(15,20.33702638248303))),
    @(16, @(@(2,16.85924304357162),@(3,10.77629454734438),@(4,6.448806841324791),@(5,11.3738818611647),@(7,27.07504590454463),@(8,13.45883590411734),@(9,20.29901313957142),@(12,9.858417560835607),@(14,16.67396544706346),@(15,20.33550852353794))),
    @(17, @(@(2,16.62918320996729),@(3,10.68598396158716),@(4,6.391471740306885),@(5,11.38701351090051),@(7,27.00372794312233),@(8,13.46895456132168),@(9,20.32978470462549),@(12,9.840961264068016),@(14,16.68931603435835),@(15,20.33605035756915))),
    @(18, @(@(2,16.49549129699225),@(3,10.63358898531044),@(4,6.358417214557016),@(5,11.3947594518248),@(7,26.96364330286156),@(8,13.47503788624505),@(9,20.34795627755056),@(12,9.831072019143791),@(14,16.69831552835907),@(15,20.33690324636186))),
    @(19, @(@(2,16.44999438012236),@(3,10.61577267553485),@(4,6.347213604729031),@(5,11.3974152282132),@(7,26.95023313303638),@(8,13.4771428002051),@(9,20.35418993273205),@(12,9.827749859593043),@(14,16.70139185822296),@(15,20.33728496311219))),
    @(20, @(@(2,16.65381589173056),@(3,10.69564459400092),@(4,6.397583423005118),@(5,11.38559565155622),@(7,27.01122330232081),@(8,13.46785015244247),@(9,20.32646008618404),@(12,9.842803922298275),@(14,16.68766432262306),@(15,20.33593665892238))),
    @(21, @(@(2,17.32202630034191),@(3,10.958536569987),@(4,6.565807148214523),@(5,11.34818954965485),@(7,27.22723394886103),@(8,13.43977619544839),@(9,20.23894485224249),@(12,9.895173794858342),@(14,16.6435744990547),@(15,20.33787364827594))),
    @(22, @(@(2,17.74579785370551),@(3,11.1260963381439),@(4,6.67478874401735),@(5,11.32546251124673),@(7,27.37647380874775),@(8,13.42377826820027),@(9,20.18597779607476),@(12,9.930698418922033),@(14,16.61627161660923),@(15,20.34395928886259))),
    @(23, @(@(2,17.52084136758495),@(3,11.03706706804784),@(4,6.616724233765894),@(5,11.33743489770324),@(7,27.29609733063331),@(8,13.43210141180827),@(9,20.21385932889708),@(12,9.911621892426099),@(14,16.63070546897704),@(15,20.3402703995722))),
    @(24, @(@(2,16.64268389360313),@(3,10.6912784964476),@(4,6.394820614095393),@(5,11.38623605408344),@(7,27.00783179112942),@(8,13.46834862724744),@(9,20.32796164938616),@(12,9.841970399306877),@(14,16.68841051871152),@(15,20.33598637571557))),
    @(25, @(@(2,15.64376959321391),@(3,10.30122706199249),@(4,6.152707932490153),@(5,11.44625235129153),@(7,26.73278178764267),@(8,13.51751126039441),@(9,20.46909482431296),@(12,9.772643898166283),@(14,16.75716589143623),@(15,20.35205097952063)))
)

foreach ($rowEntry in $data) {
    $r = $rowEntry[0]
    $cellPairs = $rowEntry[1]
    foreach ($pair in $cellPairs) {
        $c = $pair[0]
        $v = $pair[1]
        $ws.Cells.Item($r, $c).Value = $v
    }
}

Write-Host "Updated $($data.Count) rows of loading_percent values."
